$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell values in row 1 (F1, G1, H1)
$ws.Range("F1").Value = "link_organizationLink_internalRoleLinkName"
$ws.Range("G1").Value = "link_organizationLink_project_id"
$ws.Range("H1").Value = "link_organizationLink_team_id"

# Update column widths for columns F, G, H (target stored widths: 44, 34, 31)
# ColumnWidth API has a constant 0.83 offset vs the stored XML width on this sheet's font/metrics.
$ws.Columns.Item(6).ColumnWidth = 43.17
$ws.Columns.Item(7).ColumnWidth = 33.17
$ws.Columns.Item(8).ColumnWidth = 30.17
